$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999999621897643
$ws.Range("E2").Value = 0.9999999621897643

# Row 3
$ws.Range("D3").Value = 0.9999998585033933
$ws.Range("E3").Value = 0.9999998585033933

# Row 4
$ws.Range("D4").Value = 0.453001369537009
$ws.Range("E4").Value = 0.453001369537009

# Row 5
$ws.Range("D5").Value = [double]"1.211846502873738E-09"
$ws.Range("E5").Value = [double]"1.211846502873738E-09"

# Row 6
$ws.Range("D6").Value = 0.9967153217915514
$ws.Range("E6").Value = 0.9967153217915514

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.0001256165523029433
$ws.Range("E7").Value = 0.9998743834476971

# Row 8
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.05386318847245603
$ws.Range("E9").Value = 0.9461368115275439

# Row 10
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0

# Row 11
$ws.Range("D11").Value = 0.9999950180844207
$ws.Range("E11").Value = [double]"4.981915579294061E-06"
$ws.Range("F11").Value = 5.108706474304199
$ws.Range("G11").Value = 0.5

# Row 12
$ws.Range("D12").Value = 0.999999999967909
$ws.Range("E12").Value = 0.999999999967909

# Row 13
$ws.Range("D13").Value = 0.9999999999815183
$ws.Range("E13").Value = 0.9999999999815183

# Row 14
$ws.Range("D14").Value = 0.09714111744122719
$ws.Range("E14").Value = 0.09714111744122719

# Row 15
$ws.Range("D15").Value = 0.9999999999999987
$ws.Range("E15").Value = 0.9999999999999987

# Row 16
$ws.Range("D16").Value = 0.9996839453967643
$ws.Range("E16").Value = 0.9996839453967643

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = [double]"6.816241138417769E-06"
$ws.Range("E17").Value = 0.9999931837588616

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"8.949345494842818E-07"
$ws.Range("E19").Value = 0.9999991050654505

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("D21").Value = 0.9999999717960911
$ws.Range("E21").Value = [double]"2.820390887503521E-08"
$ws.Range("F21").Value = 11.71760272979736
$ws.Range("G21").Value = 0.4
